$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Mohamed Salah", "22", "5", "27"),
    @("Karim Hafez", "0", "9", "9"),
    @("Mostafa Mohamed", "8", "0", "8"),
    @("Trezeguet", "2", "1", "3"),
    @("Mohamed Elneny", "1", "0", "1"),
    @("Omar Marmoush", "0", "0", "0"),
    @("Ahmed Hegazy", "0", "0", "0"),
    @("Ahmed Elmohamady", "0", "0", "0"),
    @("Salah Basha", "0", "0", "0")
)

$row = 2
foreach ($entry in $data) {
    $rng = $ws.Range("A" + $row + ":D" + $row)
    $rng.NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]

    $rng.Style = "Normal"
    $row++
}
